$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells D1 ("A"), E1 ("B") with the same style as the
# existing header cells (bold, centered, thin border)
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D1").Value = "A"
$ws.Range("E1").Value = "B"

# Updated data values for existing columns A:C and new columns D:E
$data = @(
    @(0,   0.4, 0.5, 0.1, 0.2),
    @(0.5, 0.9, 0.4, 1.4, 2.9),
    @(0.4, 0.8, 0,   1.1, 2.3),
    @(0.8, 0.8, 0.3, 2.3, 4.6),
    @(0.6, 0.8, 0.6, 1.8, 3.6),
    @(0.2, 0.5, 0.2, 0.6, 1.1),
    @(0.9, 0.6, 0.3, 2.6, 5.2),
    @(0.3, 0.5, 0.6, 1,   2)
)

$row = 2
foreach ($rowValues in $data) {
    $col = 1
    foreach ($val in $rowValues) {
        $ws.Cells.Item($row, $col).Value = $val
        $col = $col + 1
    }
    $row = $row + 1
}
